# se agregó el reseteo de Sec cuando cambia nro de bobina
# Append two new rows (9 and 10) of bobina data to the sheet, matching
# the existing "text-typed" cell convention used throughout the sheet
# (only the "Sec" column F is stored as a true number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 -----------------------------------------------------------
# Columns A-E, G and J look numeric but must stay TEXT, same as every
# other data row already on the sheet. Writing them as a text formula
# and then collapsing the formula to its value via copy / paste-values
# keeps the literal text without leaving any numeric coercion or extra
# cell-format behind.
$ws.Range("A9").Formula = '="231"'
$ws.Range("B9").Formula = '="234"'
$ws.Range("C9").Formula = '="234"'
$ws.Range("D9").Formula = '="234"'
$ws.Range("E9").Formula = '="234"'
$ws.Range("G9").Formula = '="234"'
$ws.Range("J9").Formula = '="02"'

$ws.Range("A9:E9").Copy()
$ws.Range("A9:E9").PasteSpecial(-4163)
$ws.Range("G9").Copy()
$ws.Range("G9").PasteSpecial(-4163)
$ws.Range("J9").Copy()
$ws.Range("J9").PasteSpecial(-4163)

# "Sec" (F) is a genuine number, and the remaining columns are plain
# non-numeric text, so they can be written directly.
$ws.Range("F9").Value = 1
$ws.Range("H9").Value = "2025-02-27 17:43"
$ws.Range("I9").Value = "A"
$ws.Range("K9").Value = "COVERING"

# --- Row 10 ------------------------------------------------------------
$ws.Range("A10").Formula = '="120"'
$ws.Range("B10").Formula = '="120"'
$ws.Range("C10").Formula = '="130"'
$ws.Range("D10").Formula = '="434"'
$ws.Range("E10").Formula = '="100"'
$ws.Range("G10").Formula = '="85678"'
$ws.Range("J10").Formula = '="03"'

$ws.Range("A10:E10").Copy()
$ws.Range("A10:E10").PasteSpecial(-4163)
$ws.Range("G10").Copy()
$ws.Range("G10").PasteSpecial(-4163)
$ws.Range("J10").Copy()
$ws.Range("J10").PasteSpecial(-4163)

$ws.Range("F10").Value = 1
$ws.Range("H10").Value = "2025-03-18 14:45"
$ws.Range("I10").Value = "C"
$ws.Range("K10").Value = "L.BLANCO"
